$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.61145397019385
$ws.Range("D2").Value = 0.1636551043816326
$ws.Range("E2").Value = 1.066153845970831
$ws.Range("F2").Value = 3.24750964577936
$ws.Range("G2").Value = 0.002474875086021955
$ws.Range("L2").Value = 0.5676273488914489
$ws.Range("M2").Value = 0.4685782672052241

# Row 3
$ws.Range("B3").Value = 1.544447655151032
$ws.Range("D3").Value = 0.1570757988081937
$ws.Range("E3").Value = 0.9349628572062159
$ws.Range("F3").Value = 3.075998472660103
$ws.Range("G3").Value = 0.002486657333919957
$ws.Range("L3").Value = 0.5163001301755799
$ws.Range("M3").Value = 0.4396726270139624

# Row 4
$ws.Range("B4").Value = 1.50421495334308
$ws.Range("D4").Value = 0.1532564574921196
$ws.Range("E4").Value = 0.8542832995625531
$ws.Range("F4").Value = 2.974296035923913
$ws.Range("G4").Value = 0.002494242662180476
$ws.Range("L4").Value = 0.4850659363444549
$ws.Range("M4").Value = 0.4221872510633844

# Row 5
$ws.Range("B5").Value = 1.48804666923678
$ws.Range("D5").Value = 0.151753765122038
$ws.Range("E5").Value = 0.8213680156041505
$ws.Range("F5").Value = 2.93373080879644
$ws.Range("G5").Value = 0.002497422487713858
$ws.Range("L5").Value = 0.4724063753015457
$ws.Range("M5").Value = 0.4151269443065786

# Row 6
$ws.Range("B6").Value = 1.485375589013273
$ws.Range("D6").Value = 0.1515074397237868
$ws.Range("E6").Value = 0.8159000219439889
$ws.Range("F6").Value = 2.927047301339229
$ws.Range("G6").Value = 0.002497955869313455
$ws.Range("L6").Value = 0.4703083525190834
$ws.Range("M6").Value = 0.4139584920594146

# Row 7
$ws.Range("B7").Value = 1.503995985956351
$ws.Range("D7").Value = 0.1532359763658633
$ws.Range("E7").Value = 0.8538395524161615
$ws.Range("F7").Value = 2.973745434542934
$ws.Range("G7").Value = 0.002494285186615754
$ws.Range("L7").Value = 0.4848949298096841
$ws.Range("M7").Value = 0.4220917708379659

# Row 8
$ws.Range("B8").Value = 1.588160289157599
$ws.Range("D8").Value = 0.1613397887891921
$ws.Range("E8").Value = 1.020942235808178
$ws.Range("F8").Value = 3.187607412408397
$ws.Range("G8").Value = 0.002478865053605706
$ws.Range("L8").Value = 0.549870162785453
$ws.Range("M8").Value = 0.4585564576046508

# Row 9
$ws.Range("B9").Value = 1.760519157846886
$ws.Range("D9").Value = 0.179058423963582
$ws.Range("E9").Value = 1.34789514135224
$ws.Range("F9").Value = 3.636884731063589
$ws.Range("G9").Value = 0.002451388710044447
$ws.Range("L9").Value = 0.6796170358243216
$ws.Range("M9").Value = 0.5321982124874438

# Row 10
$ws.Range("B10").Value = 1.891765290238595
$ws.Range("D10").Value = 0.1933019100173397
$ws.Range("E10").Value = 1.588085543697872
$ws.Range("F10").Value = 3.987035239147986
$ws.Range("G10").Value = 0.002432854244336969
$ws.Range("L10").Value = 0.7765194909036381
$ws.Range("M10").Value = 0.5876815863133089

# Row 11
$ws.Range("B11").Value = 1.952511354006333
$ws.Range("D11").Value = 0.2000727703081111
$ws.Range("E11").Value = 1.69744932153651
$ws.Range("F11").Value = 4.151098936575977
$ws.Range("G11").Value = 0.002424774303784966
$ws.Range("L11").Value = 0.8209820041615785
$ws.Range("M11").Value = 0.6132399725089499

# Row 12
$ws.Range("B12").Value = 1.975666957555063
$ws.Range("D12").Value = 0.2026807178548324
$ws.Range("E12").Value = 1.738885093945981
$ws.Range("F12").Value = 4.213947174904945
$ws.Range("G12").Value = 0.002421764628644598
$ws.Range("L12").Value = 0.8378765517697673
$ws.Range("M12").Value = 0.622965544468272

# Row 13
$ws.Range("B13").Value = 1.970673163383935
$ws.Range("D13").Value = 0.2021170606477654
$ws.Range("E13").Value = 1.729960049767101
$ws.Range("F13").Value = 4.200379072855924
$ws.Range("G13").Value = 0.002422410599689255
$ws.Range("L13").Value = 0.8342354033709967
$ws.Range("M13").Value = 0.6208688488548404

# Row 14
$ws.Range("B14").Value = 1.954413309895187
$ws.Range("D14").Value = 0.2002864344747479
$ws.Range("E14").Value = 1.700857779275708
$ws.Range("F14").Value = 4.156254854602594
$ws.Range("G14").Value = 0.002424525696584745
$ws.Range("L14").Value = 0.822370758923455
$ws.Range("M14").Value = 0.6140391480133616

# Row 15
$ws.Range("B15").Value = 1.944473605527833
$ws.Range("D15").Value = 0.1991709114700484
$ws.Range("E15").Value = 1.683034909860254
$ws.Range("F15").Value = 4.129322395148449
$ws.Range("G15").Value = 0.00242582775336512
$ws.Range("L15").Value = 0.8151109015989277
$ws.Range("M15").Value = 0.6098619421705678

# Row 16
$ws.Range("B16").Value = 1.887816570913856
$ws.Range("D16").Value = 0.1928654576430233
$ws.Range("E16").Value = 1.580941024985435
$ws.Range("F16").Value = 3.976412306616453
$ws.Range("G16").Value = 0.00243338931590638
$ws.Range("L16").Value = 0.7736216624328733
$ws.Range("M16").Value = 0.5860178116599002

# Row 17
$ws.Range("B17").Value = 1.853327895180485
$ws.Range("D17").Value = 0.1890733569154577
$ws.Range("E17").Value = 1.5183406954433
$ws.Range("F17").Value = 3.883854487601866
$ws.Range("G17").Value = 0.002438117744896985
$ws.Range("L17").Value = 0.7482688864872102
$ws.Range("M17").Value = 0.5714727700869986

# Row 18
$ws.Range("B18").Value = 1.833588790781505
$ws.Range("D18").Value = 0.1869195445541152
$ws.Range("E18").Value = 1.482343565069613
$ws.Range("F18").Value = 3.831065372630576
$ws.Range("G18").Value = 0.002440870526822696
$ws.Range("L18").Value = 0.7337223382847071
$ws.Range("M18").Value = 0.5631367625549473

# Row 19
$ws.Range("B19").Value = 1.826922201010916
$ws.Range("D19").Value = 0.1861949282176454
$ws.Range("E19").Value = 1.470156859434155
$ws.Range("F19").Value = 3.813267765839328
$ws.Range("G19").Value = 0.002441808274583099
$ws.Range("L19").Value = 0.728803176342808
$ws.Range("M19").Value = 0.5603194321288214

# Row 20
$ws.Range("B20").Value = 1.856989125951941
$ws.Range("D20").Value = 0.1894741919274452
$ws.Range("E20").Value = 1.525003629266308
$ws.Range("F20").Value = 3.893660849629924
$ws.Range("G20").Value = 0.002437610972061218
$ws.Range("L20").Value = 0.7509640186186175
$ws.Range("M20").Value = 0.5730180089841213

# Row 21
$ws.Range("B21").Value = 1.959185065700694
$ws.Range("D21").Value = 0.2008229237017929
$ws.Range("E21").Value = 1.709405158821767
$ws.Range("F21").Value = 4.169195368713815
$ws.Range("G21").Value = 0.002423903087900303
$ws.Range("L21").Value = 0.8258541078284338
$ws.Range("M21").Value = 0.6160439056687892

# Row 22
$ws.Range("B22").Value = 2.026865567260757
$ws.Range("D22").Value = 0.2084972650356178
$ws.Range("E22").Value = 1.830054803305131
$ws.Range("F22").Value = 4.353491329542692
$ws.Range("G22").Value = 0.002415235518297925
$ws.Range("L22").Value = 0.8751359838992414
$ws.Range("M22").Value = 0.6444394026641049

# Row 23
$ws.Range("B23").Value = 1.990660914911587
$ws.Range("D23").Value = 0.2043770856054721
$ws.Range("E23").Value = 1.765647033560811
$ws.Range("F23").Value = 4.25473172100817
$ws.Range("G23").Value = 0.002419835080817607
$ws.Range("L23").Value = 0.8488015488837846
$ws.Range("M23").Value = 0.6292585168173019

# Row 24
$ws.Range("B24").Value = 1.85533360682831
$ws.Range("D24").Value = 0.1892928926027366
$ws.Range("E24").Value = 1.52199134040967
$ws.Range("F24").Value = 3.889226076452729
$ws.Range("G24").Value = 0.002437839977109249
$ws.Range("L24").Value = 0.7497454589473307
$ws.Range("M24").Value = 0.5723193251876779

# Row 25
$ws.Range("B25").Value = 1.713092959764083
$ws.Range("D25").Value = 0.1740580816298092
$ws.Range("E25").Value = 1.259481613704935
$ws.Range("F25").Value = 3.511955555827967
$ws.Range("G25").Value = 0.00245852929788299
$ws.Range("L25").Value = 0.6442518001229871
$ws.Range("M25").Value = 0.5120398519429656
